$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '80.985.00'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.146.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.21'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '625.81'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.88%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.283'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +23.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.579'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.145.32'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.579'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000251'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +9.63%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.27'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.726.22'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.32'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.124.85'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.154.30'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.15'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +9.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.93'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.94%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '431.42'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.99'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.07'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.15'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.16'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +7.72%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.326.57'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.02%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '75.84'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000122'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +5.67%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.98'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.25%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Cronos'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.153'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +37.61%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '557.33'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.06%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.48'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.00'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.150'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +8.54%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.68'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.406'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.92'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +8.69%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.72'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.49%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.04'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +19.07%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.99'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +10.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '160.40'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '186.61'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.32'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.89'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.771'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.51'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.03%  '
